$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 62.2
$ws.Range("K3").Value = 61.8
$ws.Range("K4").Value = 59.8
$ws.Range("K5").Value = 57.2

$ws.Range("N2").Value = 85.8724807945396
$ws.Range("N3").Value = 85.8724807945396
$ws.Range("N4").Value = 85.8724807945396
$ws.Range("N5").Value = 85.8724807945396
